# Added bar codes, serial numbers, recovery dates where needed
#
# This script fills in the Asset_Cal_Info sheet's helper lookup columns
# (J, K, O) that cross-check each Ref Des against its calibration-sheet
# name, and records two sensor bar codes (OL000333 / OL000334) that were
# previously blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# ---------------------------------------------------------------------
# 1. Newly recorded Sensor OOIBARCODE values (column E) for the VEL3D
#    (row 6) and STCENG (row 39) sensors.  These cells were previously
#    empty (just carrying a style), so clear the format before writing
#    the new text so no stray style id is left behind.
# ---------------------------------------------------------------------
$ws.Range("E6").ClearFormats()
$ws.Range("E6").Value = "OL000333"

$ws.Range("E39").ClearFormats()
$ws.Range("E39").Value = "OL000334"

# ---------------------------------------------------------------------
# 2. Column K: the distinct "calibration sheet" Ref Des names, one per
#    row 2-10, used by the cross-check MATCH() formulas below.
# ---------------------------------------------------------------------
$ws.Range("K2").Value  = "CP02PMUI-SBS01-00-STCENG000"
$ws.Range("K3").Value  = "CP02PMUI-SBS01-01-MOPAK0000"
$ws.Range("K4").Value  = "CP02PMUI-RII01-02-ADCPTG010"
$ws.Range("K5").Value  = "CP02PMUI-WFP01-00-WFPENG000"
$ws.Range("K6").Value  = "CP02PMUI-WFP01-01-VEL3DK000"
$ws.Range("K7").Value  = "CP02PMUI-WFP01-02-DOFSTK000"
$ws.Range("K8").Value  = "CP02PMUI-WFP01-03-CTDPFK000"
$ws.Range("K9").Value  = "CP02PMUI-WFP01-04-FLORTK000"
$ws.Range("K10").Value = "CP02PMUI-WFP01-05-PARADK000"

# ---------------------------------------------------------------------
# 3. Column J: for every data row (2-39), does the Ref Des in column A
#    show up somewhere in the new K list?  Row 2 is the standalone
#    "master" formula; rows 3-39 are entered as one shared formula.
# ---------------------------------------------------------------------
$ws.Range("J2").Formula   = "=MATCH(A2,K:K,0)"
$ws.Range("J3:J39").Formula = "=MATCH(A3,K:K,0)"

# ---------------------------------------------------------------------
# 4. Column O: the inverse check - where does each K entry show up in
#    column A?  Again row 2 is standalone, rows 3-10 share a formula.
# ---------------------------------------------------------------------
$ws.Range("O2").Formula    = "=MATCH(K2,A:A,0)"
$ws.Range("O3:O10").Formula = "=MATCH(K3,A:A,0)"

# Put the cursor where it ended up in the authored workbook and make
# sure everything above is recalculated.
$ws.Range("N15").Select()
$excel.Calculate()
